# Generate Report for Handoff
# bf804513-... has just been handed off: status flips from "Ready for handoff"
# to "In Translation" and its handoff timestamps are refreshed.
# fae22e52-... was the next file in the queue and is now "Ready for handoff".

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E6").Value = "In Translation"
$ws.Range("F6").Value = "In Translation"
$ws.Range("G6").Value = "2016-10-20 00:04:56"
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C6").Value = "In Translation"
$ws.Range("H6").Value = "2016-10-20 00:04:45"
$ws.Range("C7").Value = "Ready for handoff"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C6").Value = "In Translation"
$ws.Range("H6").Value = "2016-10-20 00:04:56"
$ws.Range("C7").Value = "Ready for handoff"
